$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44174
$ws.Range("J2").Value = 2800
$ws.Range("L2").Value = 1250
$ws.Range("M2").Value = 1221
$ws.Range("P2").Value = 1221
$ws.Range("D3").Value = 44174
$ws.Range("J3").Value = 1300
$ws.Range("D4").Value = 44179
$ws.Range("J4").Value = 980
$ws.Range("K4").Value = 1200
$ws.Range("L4").Value = 1200
$ws.Range("M4").Value = 1200
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 1200
$ws.Range("D5").Value = 44176
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 2500
$ws.Range("K5").Value = 1200
$ws.Range("L5").Value = 1300
$ws.Range("M5").Value = 1256
$ws.Range("P5").Value = 1256
$ws.Range("D6").Value = 44176
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 1500
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = 1000
$ws.Range("P6").Value = 1000
$ws.Range("D7").Value = 44175
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 1500
$ws.Range("K7").Value = 1300
$ws.Range("L7").Value = 1300
$ws.Range("M7").Value = 1300
$ws.Range("P7").Value = 1300
$ws.Range("D8").Value = 44175
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 1450
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 1000
$ws.Range("M8").Value = 1000
$ws.Range("P8").Value = 1000
$ws.Range("D9").Value = 44168
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 1200
$ws.Range("K9").Value = 1300
$ws.Range("L9").Value = 1300
$ws.Range("M9").Value = 1300
$ws.Range("P9").Value = 1300
$ws.Range("D10").Value = 44168
$ws.Range("I10").Value = "Segunda"
$ws.Range("J10").Value = 850
$ws.Range("K10").Value = 1000
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = 1000
$ws.Range("O10").Value = "Provincia de Quillota"
$ws.Range("P10").Value = 1000
$ws.Range("D11").Value = 44161
$ws.Range("J11").Value = 1600
$ws.Range("D12").Value = 44161
$ws.Range("J12").Value = 1850
$ws.Range("D13").Value = 44172
$ws.Range("J13").Value = 600
$ws.Range("D14").Value = 44172
$ws.Range("J14").Value = 550
$ws.Range("D15").Value = 44169
$ws.Range("J15").Value = 950
$ws.Range("K15").Value = 1300
$ws.Range("M15").Value = 1300
$ws.Range("P15").Value = 1300
$ws.Range("D16").Value = 44169
$ws.Range("J16").Value = 800
$ws.Range("D17").Value = 44181
$ws.Range("J17").Value = 1000
$ws.Range("D18").Value = 44181
$ws.Range("J18").Value = 900
$ws.Range("K18").Value = 900
$ws.Range("L18").Value = 900
$ws.Range("M18").Value = 900
$ws.Range("P18").Value = 900
$ws.Range("D19").Value = 44162
$ws.Range("J19").Value = 1200
$ws.Range("D20").Value = 44162
$ws.Range("J20").Value = 800
$ws.Range("D21").Value = 44167
$ws.Range("J21").Value = 1430
$ws.Range("K21").Value = 1200
$ws.Range("M21").Value = 1248
$ws.Range("P21").Value = 1248
$ws.Range("D22").Value = 44167
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = 1000
$ws.Range("P22").Value = 1000
$ws.Range("D25").Value = 44159
$ws.Range("J25").Value = 1100
$ws.Range("D26").Value = 44159
$ws.Range("D27").Value = 44165
$ws.Range("J27").Value = 720
$ws.Range("L27").Value = 1200
$ws.Range("M27").Value = 1200
$ws.Range("P27").Value = 1200
$ws.Range("D28").Value = 44165
$ws.Range("J28").Value = 750